$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the KV (Krankenversicherung) related data block (rows 39-49)
$ws.Range("B39").Value = "nein"

$ws.Range("B40").Value = $null
$ws.Range("B41").Value = $null
$ws.Range("B42").Value = $null

$ws.Range("B43").Value = "ja"

$ws.Range("A44").Value = "AG-Krankenversicherungsbeitrag in Prozent"
$ws.Range("B44").Value = 7.3

$ws.Range("A45").Value = "AN-Krankenversicherungsbeitrag in Prozent"
$ws.Range("B45").Value = 7.3

$ws.Range("A46").Value = "Beitragsbemessungsgrenze Krankenversicherung Ost"
$ws.Range("B46").Value = 72000

$ws.Range("A47").Value = "Beitragsbemessungsgrenze Krankenversicherung West"
$ws.Range("B47").Value = 68000

$ws.Range("A48").Value = "Mitglied gesetzliche Krankenkasse (vollständiger Name)"

$ws.Range("A49").Value = "Mitglied gesetzliche Krankenkasse (Abkürzung)"

# Update the active sheet view / selection
$ws.Activate()
$ws.Range("B41").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
